# Adds a "FOOT" (footprint) sprite block (columns BV:CF) to the sprite sheet,
# and performs minor view/range housekeeping, per commit "add footprint and fix range graphics".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Reference cells whose formatting we reuse for the new sprite block
$srcHeader = $ws.Range("B2")   # section header style (grey fill, centered)
$srcK = $ws.Range("B4")        # "K" (black/outline) pixel style
$srcD = $ws.Range("AX4")       # "D" (pad/shadow) pixel style

# --- New column block BV:CF gets the same narrow pixel-art column width as the other sprites ---
$ws.Range("BV1:CF1").ColumnWidth = $ws.Range("B1").ColumnWidth

# --- Row 2: new "FOOT" section header, merged across BV2:CF2 ---
$hdrRange = $ws.Range("BV2:CF2")
$srcHeader.Copy()
$hdrRange.PasteSpecial($xlPasteFormats)
$ws.Range("BV2").Value = "FOOT"
$hdrRange.Merge()

# --- Rows 4-14: draw the FOOT sprite (11x11 pixel grid) in BV:CF ---
# Base fill: every pixel starts as "K"
$kRange = $ws.Range("BV4:CF14")
$srcK.Copy()
$kRange.PasteSpecial($xlPasteFormats)
$kRange.Value = "K"

# Overlay the footprint pad/toes as "D" pixels
$srcD.Copy()
$ws.Range("CA5:CA5").PasteSpecial($xlPasteFormats)
$ws.Range("CA5:CA5").Value = "D"
$ws.Range("BW6:BX7").PasteSpecial($xlPasteFormats)
$ws.Range("BW6:BX7").Value = "D"
$ws.Range("BZ6:CB7").PasteSpecial($xlPasteFormats)
$ws.Range("BZ6:CB7").Value = "D"
$ws.Range("CD6:CE7").PasteSpecial($xlPasteFormats)
$ws.Range("CD6:CE7").Value = "D"
$ws.Range("BZ9:CB9").PasteSpecial($xlPasteFormats)
$ws.Range("BZ9:CB9").Value = "D"
$ws.Range("BY10:CC10").PasteSpecial($xlPasteFormats)
$ws.Range("BY10:CC10").Value = "D"
$ws.Range("BX11:CD13").PasteSpecial($xlPasteFormats)
$ws.Range("BX11:CD13").Value = "D"
$ws.Range("BY14:CC14").PasteSpecial($xlPasteFormats)
$ws.Range("BY14:CC14").Value = "D"

# --- Column CG acts as the separator ( " ) to the right of every sprite row, like columns M, Y, AK, AW, BI, BU ---
$quote = [char]34
$ws.Range("CG4:CG14").Value = $quote
$ws.Range("CG16:CG26").Value = $quote

# --- View / selection housekeeping to match the saved workbook state ---
$win = $excel.ActiveWindow
$win.Zoom = 80
$win.ScrollRow = 3
$win.ScrollColumn = 15
$ws.Range("BX29").Select() | Out-Null

